$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "teste"
$ws.Range("D2").Value = "zap"
$ws.Range("C3").Value = "teste"
$ws.Range("D3").Value = "zap"

$ws.Range("A4").Value = "Cristian"
$ws.Range("B4").Value = 5543996777718
$ws.Range("C4").Value = "teste"
$ws.Range("D4").Value = "zap"
